# Remove the trailing "Ver no Jupiter..." / copyright footer block, along
# with the blank paragraph that separated it from the "Requisitos" text,
# mirroring the site rebuild that dropped this boilerplate from the page.

$d = $word.ActiveDocument

$targetText = "Ver no Jupiter Salvar em pdf Salvar em docx"

$foundIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$targetText*") {
        $foundIndex = $i
        break
    }
}

if ($foundIndex -gt 0) {
    # Also remove the blank paragraph right before it and the copyright
    # paragraph right after it (three whole paragraphs in total).
    $startPara = $d.Paragraphs.Item($foundIndex - 1)
    $endPara = $d.Paragraphs.Item($foundIndex + 1)

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
